$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.706.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "'1.887.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'247.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.4732"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "'0.06527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "'22.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "'0.07798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "'1.891.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "'96.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "'0.7349"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").Value = "'283.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.80%  "
$ws.Range("D17").Value = "'30.699.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "'13.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "'0.000007528"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'2.141.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").Value = "'5.314"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'6.246"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'9.220"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'164.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "'18.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "'1.920"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "'1.340"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").Value = "'0.09729"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "'1.491"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "'4.298"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'4.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "'0.04867"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "'1.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").Value = "'0.6978"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "'0.01892"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").Value = "'6.369"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "'76.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("D42").Value = "'2.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'0.4251"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'0.8360"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "'101.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("D47").Value = "'9.527"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D48").Value = "'35.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "'7.025"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "'915.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "'0.05756"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.19%  "
